# Adds the "Python + Windows Script Host (WSH)" Q&A section to Sheet1,
# appended below the existing content (rows 126-153, column B).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B126").Value = "Kết hợp Python với Windows Script Host (WSH) để tự động hóa trình duyệt IE có thể là một gợi ý thú vị. Tuy nhiên, chúng ta sẽ sử dụng Python để viết mã tự động hóa và sử dụng WSH để chạy mã Python. Dưới đây là các bước mà bạn có thể thực hiện:"
$ws.Range("B128").Value = "Viết mã tự động hóa bằng Python: Đầu tiên, bạn cần viết mã Python để tự động hóa trình duyệt IE bằng cách sử dụng thư viện Selenium WebDriver. Dưới đây là một ví dụ về việc sử dụng Python để mở trình duyệt IE và mở một trang web:"
$ws.Range("B129").Value = "pythonCopy code"
$ws.Range("B130").Value = "from selenium import webdriver"
$ws.Range("B131").Value = "from selenium.webdriver.ie.options import Options as IEOptions"
$ws.Range("B133").Value = "options = IEOptions()"
$ws.Range("B134").Value = "options.ignore_protected_mode_settings = True"
$ws.Range("B135").Value = "driver = webdriver.Ie(executable_path='path/to/IEDriverServer.exe', options=options)"
$ws.Range("B137").Value = "driver.get('https://www.example.com')"
$ws.Range("B139").Value = "# Thực hiện các tác vụ tự động hóa khác tại đây"
$ws.Range("B141").Value = "driver.quit()"
$ws.Range("B142").Value = "Tạo tệp mã Python: Bạn cần tạo một tệp văn bản (ví dụ: automate_ie.py) và chèn mã tự động hóa viết bằng Python vào tệp đó."
$ws.Range("B143").Value = "Tạo mã WSH để chạy tệp Python: Sau khi bạn đã có tệp mã Python, bạn có thể viết mã WSH để chạy tệp Python đó. Dưới đây là ví dụ về cách bạn có thể thực hiện điều này:"
$ws.Range("B144").Value = "vbscriptCopy code"
$ws.Range("B145").Value = "Dim objShell"
$ws.Range("B146").Value = "Set objShell = CreateObject(`"WScript.Shell`")"
$ws.Range("B147").Value = "objShell.Run `"python path\to\automate_ie.py`", 1, True"
$ws.Range("B148").Value = "Set objShell = Nothing"
$ws.Range("B149").Value = "Trong mã trên, path\to\automate_ie.py là đường dẫn đến tệp mã Python bạn đã viết."
$ws.Range("B150").Value = "Lưu tệp mã WSH: Lưu tệp mã WSH với phần mở rộng .vbs (ví dụ: run_automate_ie.vbs)."
$ws.Range("B151").Value = "Chạy mã WSH: Bạn có thể chạy tệp mã WSH bằng cách nhấp đôi vào nó, hoặc bạn cũng có thể chạy nó từ dòng lệnh bằng cách gõ cscript run_automate_ie.vbs."
$ws.Range("B153").Value = "Lưu ý rằng việc kết hợp Python với WSH có thể không phải là cách tiếp cận thông thường, và nó có thể gây ra một số khó khăn trong việc quản lý và điều hướng mã. Tuy nhiên, nếu bạn muốn sử dụng Python để viết mã tự động hóa trong môi trường Windows Script Host, đây là một trong những cách bạn có thể thử."

# Matches the author's final selection/view position recorded in the diff.
$ws.Range("D125").Select() | Out-Null

